$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (89) down to the new row (90)
$ws.Range("A89:V89").Copy()
$ws.Range("A90:V90").PasteSpecial(-4122)

$ws.Range("A90").Value = 89
$ws.Range("B90").Value = "bosnia-and-herzegovina"
$ws.Range("C90").Value = "premijer-liga-bih"
$ws.Range("D90").Value = "2023-2024"
$ws.Range("E90").Value = 45261.75
$ws.Range("F90").Value = "Velez Mostar"
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = "FK Sarajevo"
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 2.27
$ws.Range("K90").Value = "30/11/2023 06:12"
$ws.Range("L90").Value = 2.21
$ws.Range("M90").Value = "01/12/2023 17:50"
$ws.Range("N90").Value = 2.99
$ws.Range("O90").Value = "30/11/2023 06:12"
$ws.Range("P90").Value = 3.14
$ws.Range("Q90").Value = "01/12/2023 17:50"
$ws.Range("R90").Value = 3.01
$ws.Range("S90").Value = "30/11/2023 06:12"
$ws.Range("T90").Value = 3.34
$ws.Range("U90").Value = "01/12/2023 17:50"
$ws.Range("V90").Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/velez-mostar-fk-sarajevo/vgYBD2iq/"
